$wb = $excel.ActiveWorkbook

# ---- Portal sheet: no longer the active tab; selection becomes the used range A1:C9 ----
$portal = $wb.Worksheets.Item("Portal")
$portal.Range("A1:C9").Select() | Out-Null

# ---- New sheet: "Field Scouting" (inserted after Portal) ----
$fieldScouting = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$fieldScouting.Name = "Field Scouting"

$fieldScouting.Columns.Item(1).ColumnWidth = 26.1640625

$fieldScouting.Range("A1").Value = "Feature"
$fieldScouting.Range("B1").Value = "Works"
$fieldScouting.Range("C1").Value = "Notes"

$fieldScouting.Range("A2").Value = "Field Schedule"
$portal.Range("B2").Copy($fieldScouting.Range("B2")) | Out-Null
$fieldScouting.Range("B2").Value = 45317

$fieldScouting.Range("A3").Value = "Auto Team Select By Match"
$portal.Range("B2").Copy($fieldScouting.Range("B3")) | Out-Null
$fieldScouting.Range("B3").Value = 45317

$fieldScouting.Range("A4").Value = "Save"
$portal.Range("B2").Copy($fieldScouting.Range("B4")) | Out-Null
$fieldScouting.Range("B4").Value = 45317

$portal.Range("B2").Copy($fieldScouting.Range("B5")) | Out-Null
$fieldScouting.Range("B5").ClearContents() | Out-Null

$portal.Range("B2").Copy($fieldScouting.Range("B6")) | Out-Null
$fieldScouting.Range("B6").ClearContents() | Out-Null

$portal.Range("B2").Copy($fieldScouting.Range("B7")) | Out-Null
$fieldScouting.Range("B7").ClearContents() | Out-Null

$portal.Range("B2").Copy($fieldScouting.Range("B9")) | Out-Null
$fieldScouting.Range("B9").ClearContents() | Out-Null

$fieldScouting.Range("A1:C4").Select() | Out-Null

# ---- New sheet: "Scout Field Results" (inserted after Field Scouting, becomes active tab) ----
$scoutFieldResults = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$scoutFieldResults.Name = "Scout Field Results"

$scoutFieldResults.Range("A1").Value = "Feature"
$scoutFieldResults.Range("B1").Value = "Works"
$scoutFieldResults.Range("C1").Value = "Notes"

$scoutFieldResults.Range("A2:XFD4").Select() | Out-Null

Write-Host "Edit complete."
